# Update "想去人数" (want-to-go count) figures that changed between crawls.
# Both the "展览" sheet and the "全部类型" sheet (which aggregates all rows)
# need the same four rows updated.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F4").Value = 253
    $ws.Range("F16").Value = 5394
    $ws.Range("F26").Value = 5
    $ws.Range("F29").Value = 67
}

$wb.Save()
